$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (A2:Z2) - automatic electricity price update
$ws.Range("A2").Value = 45948
$ws.Range("B2").Value = 125.52
$ws.Range("C2").Value = 112.9
$ws.Range("D2").Value = 100.88
$ws.Range("E2").Value = 101.35
$ws.Range("F2").Value = 102.34
$ws.Range("G2").Value = 97.45
$ws.Range("H2").Value = 103.83
$ws.Range("I2").Value = 114.9
$ws.Range("J2").Value = 118.52
$ws.Range("K2").Value = 104.94
$ws.Range("L2").Value = 75.37
$ws.Range("M2").Value = 29.27
$ws.Range("N2").Value = 6.74
$ws.Range("O2").Value = 4.66
$ws.Range("P2").Value = 4.31
$ws.Range("Q2").Value = 1.89
$ws.Range("R2").Value = 5.91
$ws.Range("S2").Value = 47.97
$ws.Range("T2").Value = 96.20999999999999
$ws.Range("U2").Value = 117.97
$ws.Range("V2").Value = 128.65
$ws.Range("W2").Value = 120.73
$ws.Range("X2").Value = 113.17
$ws.Range("Y2").Value = 109.17
$ws.Range("Z2").Value = 81.03

# AA2 unchanged (20h-24h)

$ws.Range("AB2").Value = 117.93

# AC2 unchanged (20h-22h)

$ws.Range("AD2").Value = 124.69
$ws.Range("AE2").Value = "0h-2h"
$ws.Range("AF2").Value = 119.21
$ws.Range("AG2").Value = "10h-17h"
